$wb = $excel.ActiveWorkbook

$wsMeta = $wb.Worksheets.Item("Metadata")
$wsElem = $wb.Worksheets.Item("Elements")

# --- Metadata sheet updates ---
# Title
$wsMeta.Range("B5").Value = "Range with UCUM or EDQM codes if code is used"
# Date
$wsMeta.Range("B8").Value = "2025-08-13T14:10:49+00:00"
# Description
$wsMeta.Range("B12").Value = "Range with low and high unit UCUM or EDQM codes if code is used"

# --- Elements sheet updates ---
# Row 2 (Range element): Short / Definition / Comments - drop "with UCUM or EDQM unit" wording
$wsElem.Range("L2").Value = "Set of values bounded by precise low and high fixed quantity (no comparator)"
$wsElem.Range("M2").Value = "A set of ordered Quantities defined by a precise low and high limit defined by a fixed quantity (no comparator)"
$wsElem.Range("N2").Value = "The stated low and high value are assumed to have arbitrarily high precision when it comes to determining which values are in the range. I.e. 1.99 is not in the range 2 -> 3. Low and high limit are precisely defined, no element 'comparator' in the simpleQuantity defining each bound. The limits are defined by a fixed quantity (no comparator)."

# Row 5 (Range.low): Short / Definition rewritten, Requirements cleared
$wsElem.Range("L5").Value = "Low limit"
$wsElem.Range("M5").Value = "The low limit. The boundary is inclusive."
$wsElem.Range("O5").ClearContents()

# Row 6 (Range.high): Short / Definition rewritten, Requirements cleared
$wsElem.Range("L6").Value = "High limit"
$wsElem.Range("M6").Value = "The high limit. The boundary is inclusive."
$wsElem.Range("O6").ClearContents()

# Column O ("Requirements") no longer holds any long text once O5/O6 are
# cleared, so its best-fit width shrinks. Approximate Excel's recalculated
# best-fit width for the narrower column.
$wsElem.Columns.Item(15).ColumnWidth = 12.59
